# Update "paises.xlsx" (COVID-19 country stats) to a newer data snapshot.
# - Refresh the "last updated" timestamp.
# - Update case/death statistics for a handful of countries.
# - Re-sort the country table by "Casos totales" (column B) descending, since
#   the sheet is expected to stay ranked by total cases and a couple of
#   countries changed rank because of the new numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Refresh the timestamp banner in A1.
$ws.Range("A1").Value = "Datos actualizados a 25 de Marzo de 2020 a las 08:46"

# 2) Helper: find the row whose column-A country name matches, then write
#    the new Casos totales / Nuevos casos / Casos activos / Recuperados /
#    Casos criticos / Muertes hoy / Muertes values (columns B..H).
#    (Positional parameters only -- named parameter binding isn't reliable
#    for script-defined functions in this host.)
function Set-CountryRow($Sheet, $Country, $Values) {
    $row = $null
    for ($r = 4; $r -le 200; $r++) {
        if ($Sheet.Cells.Item($r, 1).Value() -eq $Country) {
            $row = $r
            break
        }
    }

    if ($null -eq $row) {
        return
    }

    for ($i = 0; $i -lt $Values.Length; $i++) {
        $Sheet.Cells.Item($row, 2 + $i).Value = $Values[$i]
    }
}

Set-CountryRow $ws "Suiza"      @(9891, 14, 131, 9628, 141, 10, 132)
Set-CountryRow $ws "Austria"    @(5315, 32, 9,   5276, 26,  2,  30)
Set-CountryRow $ws "Australia"  @(2423, 106, 118, 2297, 11,  0,  8)
Set-CountryRow $ws "Malasia"    @(1796, 172, 183, 1596, 64,  1,  17)
Set-CountryRow $ws "Tailandia"  @(934,  107, 70,  860,  11,  0,  4)
Set-CountryRow $ws "Serbia"     @(303,  0,   15,  284,  21,  1,  4)
Set-CountryRow $ws "Taiwan"     @(235,  19,  29,  204,  0,   0,  2)
Set-CountryRow $ws "Lituania"   @(255,  46,  1,   252,  1,   0,  2)
Set-CountryRow $ws "Ucrania"    @(113,  11,  1,   109,  0,   0,  3)
Set-CountryRow $ws "Oman"       @(99,   15,  17,  82,   0,   0,  0)
Set-CountryRow $ws "Afganistan" @(74,   0,   1,   71,   0,   1,  2)
Set-CountryRow $ws "Georgia"    @(73,   3,   10,  63,   1,   0,  0)
Set-CountryRow $ws "Nigeria"    @(46,   2,   2,   43,   0,   0,  1)

# 3) Re-sort the data rows (A4:H200) by Casos totales (column B) descending,
#    keeping countries with equal totals in their existing relative order.
$dataRange = $ws.Range("A4:H200")
$keyRange = $ws.Range("B4:B200")
$dataRange.Sort($keyRange, 2)
